$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the prompt text stored in E5 (shared string) with the new task description.
$newText = @"
Provide a list of requested number of items.
任务目标: 设计不同人物对相关主题的看法，先用中高级日语回答，然后提供中文翻译。然后对重点单词或语法做说明。
Number of items: 10
Example output item:
'''
item x: {日语句子}/{中文翻译}/{重点单词或语法说明}
'''
Reply in the following format:
    - item 1
    - item 2
    - item 3
"@

$ws.Range("E5").Value = $newText

# Move the active selection from E6 to H5 (view stays scrolled at A5).
$ws.Range("H5").Select() | Out-Null
